# Add a new "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/bordered/centered style
# used by the other header cells, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data value for row 2 (plain number, no special formatting).
$ws.Range("H2").Value = 1
